# [CSTPER-451] Import via XLS should offer stop in workspace, send to workflow and skip workflow
#
# Add two new metadata columns ("dc.type" and "dc.date.issued") to the
# "Main" sheet of the bulk-import example workbook, with sample values
# for the two existing data rows, and make "Main" the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# New header cells
$ws.Range("D1").Value = "dc.type"
$ws.Range("E1").Value = "dc.date.issued"

# New sample data
$ws.Range("D2").Value = "Article"
$ws.Range("D3").Value = "Book"

$ws.Range("E2").Value = 43831
$ws.Range("E3").Value = 44298
$ws.Range("E2:E3").NumberFormat = "yyyy\-mm\-dd"

# "Main" becomes the selected/active sheet (previously it was
# "dc.contributor.author")
[void]$ws.Range("E4").Select()
$ws.Activate()
